$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add a new table column ("Column1") at the end of Table1
$tbl = $ws.ListObjects.Item("Table1")
$newCol = $tbl.ListColumns.Add()
$newCol.Range.Cells.Item(1,1).Value = "Column1"

# 2. Enter the weakness note about Sharding (row 9, Weaknesses column = G)
$ws.Range("G9").Value = "Multiple shards must work together, increasing security risks, If one shard is compromised, the whole network is threatened "

# 3. Fill the remaining data cells (B2:H10, minus G9) with "N/A"
$ws.Range("B2:H10").Value = "N/A"
$ws.Range("G9").Value = "Multiple shards must work together, increasing security risks, If one shard is compromised, the whole network is threatened "

# 4. Remove the now-unused Percent cell style / number format from column E,
#    and wrap the long weakness text so it is fully visible
$ws.Range("E2:E10").Style = "Normal"
$wb.Styles.Item("Percent").Delete()
$ws.Range("G9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 170

# 5. Update the active selection to match the end state
$ws.Range("J9").Select()
